$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 599
$ws1.Range("G4").Value = 45
$ws1.Range("F5").Value = 540
$ws1.Range("F6").Value = 301
$ws1.Range("F7").Value = 2718
$ws1.Range("F9").Value = 7575
$ws1.Range("F11").Value = 462
$ws1.Range("F12").Value = 30
$ws1.Range("F13").Value = 271

# Sheet "全部类型" (all types listing, mirrors exhibition rows + others)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 599
$ws4.Range("G4").Value = 45
$ws4.Range("F5").Value = 540
$ws4.Range("F6").Value = 301
$ws4.Range("F9").Value = 2718
$ws4.Range("F11").Value = 7576
$ws4.Range("F13").Value = 462
$ws4.Range("F14").Value = 30
$ws4.Range("F17").Value = 271
